$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-09-22T01:31:54.822615", "TRADING_ATTEMPT", "SOL", "UNKNOWN", 236.6486785425042, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 1/3"),
    @("2025-09-22T01:31:56.476164", "POSITION_OPENED", "SOL", "UNKNOWN", 236.6486785425042, 2400, 20, 0.09115590152285132, $null, $null, "SUCCESS", $null),
    @("2025-09-22T01:31:56.494749", "TRADING_ATTEMPT", "BTC", "UNKNOWN", 115304.4799940704, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 2/3"),
    @("2025-09-22T01:31:58.040397", "POSITION_OPENED", "BTC", "UNKNOWN", 115304.4799940704, 4800, 40, 0.2552237794056823, $null, $null, "SUCCESS", $null),
    @("2025-09-22T01:31:58.059394", "TRADING_ATTEMPT", "ENA", "UNKNOWN", 0.6448516400994989, $null, $null, $null, $null, $null, "ATTEMPT", "Attempting trade 3/3"),
    @("2025-09-22T01:31:59.674594", "POSITION_OPENED", "ENA", "UNKNOWN", 0.6448516400994989, 2400, 10, 1.301833833851505, $null, $null, "SUCCESS", $null)
)

$startRow = 26
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $val = $values[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}
